$wb = $excel.ActiveWorkbook
$wsSearch = $wb.Worksheets.Item("NaukriSearch")
$wsLog = $wb.Worksheets.Item("Sheet1")

# --- 1) Append a new row (12) onto the history log in Sheet1, cloning the
#        formatting (borders / number formats / wrap text / hyperlink style)
#        of the existing result row on the NaukriSearch tab before that row's
#        own content gets replaced below. ---
$wsSearch.Range("A2:K2").Copy($wsLog.Range("A12:K12"))

$wsLog.Range("A12").Value = 1
$wsLog.Range("B12").Value = "Workday Support Analyst"
$wsLog.Range("C12").Value = "Workday Support Analyst"
$wsLog.Range("D12").Value = "Workday Support Analyst"
$wsLog.Range("E12").Value = "4-10"
$wsLog.Range("F12").Value = 15
$wsLog.Range("G12").Value = "Bangalore"
$wsLog.Range("H12").Value = "1 Month"
$wsLog.Range("I12").Value = 25
$wsLog.Range("J12").Value = "15-20"
$wsLog.Range("K12").Value = "sivaprakasam.sundaram@cai.io"
$wsLog.Range("K12").Hyperlinks.Item(1).Address = "mailto:sivaprakasam.sundaram@cai.io"

$wsLog.Rows.Item(12).RowHeight = 43.2

# --- 2) Replace the "latest result" row shown on the NaukriSearch tab with the
#        newest job listing that was just found. ---
$wsSearch.Range("B2").Value = "Sitecore Developer"
$wsSearch.Range("C2").Value = "Sitecore"
$wsSearch.Range("D2").Value = "Sitecore"
$wsSearch.Range("E2").Value = "6-10"
$wsSearch.Range("H2").Value = "15 Days"
$wsSearch.Range("J2").Value = "7-24"
$wsSearch.Range("K2").Value = "sivaprakasam.sundaram@cai.io"
$wsSearch.Range("K2").Hyperlinks.Delete()
$wsSearch.Range("K2").Hyperlinks.Add($wsSearch.Range("K2"), "mailto:sivaprakasam.sundaram@cai.io") | Out-Null
$wsSearch.Range("K2").Borders.LineStyle = 1

$wsSearch.Range("F2").Select() | Out-Null
